{"js": "// Office.js (Word JavaScript API) edit script.\n// This is the body of `async (context) => { ... }`.\n//\n// The document is a \"two-digit division\" worksheet: a title paragraph with a\n// date, followed by a single 20-row x 5-col table whose 5 \"problem\" rows\n// (0, 4, 8, 12, 16) hold \"A\u00f7B=\" text runs. We update the date and 25 of the\n// division problems to new values.\n//\n// NOTE: several new values happen to equal OTHER cells' *old* values (e.g.\n// cell (0,2) \"35\u00f73=\" becomes \"11\u00f76=\", while cell (12,1) already holds\n// \"11\u00f76=\" and itself becomes \"49\u00f75=\"). A blind document-wide text\n// search-and-replace can therefore double-hit a freshly written cell. To\n// stay correct regardless of execution order we address each table cell by\n// its stable (row, column) position instead of searching by text.\n\n// Row/column position (0-based) -> new text for each division-problem cell,\n// addressed directly so textual collisions between old/new values can't\n// cause a mis-replace.\nconst cellUpdates = [\n  [0, 0, \"64\u00f79=\"],\n  [0, 1, \"98\u00f76=\"],\n  [0, 2, \"11\u00f76=\"],\n  [0, 3, \"39\u00f77=\"],\n  [0, 4, \"24\u00f72=\"],\n  [4, 0, \"52\u00f78=\"],\n  [4, 1, \"50\u00f79=\"],\n  [4, 2, \"64\u00f73=\"],\n  [4, 3, \"29\u00f73=\"],\n  [4, 4, \"88\u00f78=\"],\n  [8, 0, \"48\u00f76=\"],\n  [8, 1, \"70\u00f78=\"],\n  [8, 2, \"71\u00f72=\"],\n  [8, 3, \"53\u00f78=\"],\n  [8, 4, \"90\u00f76=\"],\n  [12, 0, \"87\u00f76=\"],\n  [12, 1, \"49\u00f75=\"],\n  [12, 2, \"67\u00f74=\"],\n  [12, 3, \"74\u00f78=\"],\n  [12, 4, \"56\u00f76=\"],\n  [16, 0, \"18\u00f73=\"],\n  [16, 1, \"86\u00f79=\"],\n  [16, 2, \"99\u00f78=\"],\n  [16, 3, \"58\u00f79=\"],\n  [16, 4, \"92\u00f79=\"],\n];\n\n// Update the title date paragraph (the document's first paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2023-12-01 Friday\", \"Replace\");\n\n// Update each division-problem table cell by its fixed (row, column).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nfor (const [row, col, newText] of cellUpdates) {\n  table.getCell(row, col).insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The document is a \"two-digit division\" worksheet: a title paragraph with a\n# date, followed by a single 20-row x 5-col table whose 5 \"problem\" rows\n# (1, 5, 9, 13, 17 in 1-based COM numbering) hold \"A\u00f7B=\" text runs. We update\n# the date and 25 of the division problems to new values.\n#\n# NOTE: several new values happen to equal OTHER cells' *old* values (e.g.\n# cell (1,3) \"35\u00f73=\" becomes \"11\u00f76=\", while cell (13,2) already holds\n# \"11\u00f76=\" and itself becomes \"49\u00f75=\"). A blind Find/Replace across the whole\n# document could therefore double-hit a freshly written cell. To stay\n# correct regardless of execution order we address each table cell directly\n# by its stable (row, column) position via Table.Cell(row, col) instead of\n# searching by text.\n\n$d = $word.ActiveDocument\n\n# Update the title date paragraph (the document's first paragraph).\n$d.Paragraphs.Item(1).Range.Text = \"2023-12-01 Friday\"\n\n# Update each division-problem table cell by its fixed (row, column).\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"64\u00f79=\"\n$t.Cell(1, 2).Range.Text = \"98\u00f76=\"\n$t.Cell(1, 3).Range.Text = \"11\u00f76=\"\n$t.Cell(1, 4).Range.Text = \"39\u00f77=\"\n$t.Cell(1, 5).Range.Text = \"24\u00f72=\"\n\n$t.Cell(5, 1).Range.Text = \"52\u00f78=\"\n$t.Cell(5, 2).Range.Text = \"50\u00f79=\"\n$t.Cell(5, 3).Range.Text = \"64\u00f73=\"\n$t.Cell(5, 4).Range.Text = \"29\u00f73=\"\n$t.Cell(5, 5).Range.Text = \"88\u00f78=\"\n\n$t.Cell(9, 1).Range.Text = \"48\u00f76=\"\n$t.Cell(9, 2).Range.Text = \"70\u00f78=\"\n$t.Cell(9, 3).Range.Text = \"71\u00f72=\"\n$t.Cell(9, 4).Range.Text = \"53\u00f78=\"\n$t.Cell(9, 5).Range.Text = \"90\u00f76=\"\n\n$t.Cell(13, 1).Range.Text = \"87\u00f76=\"\n$t.Cell(13, 2).Range.Text = \"49\u00f75=\"\n$t.Cell(13, 3).Range.Text = \"67\u00f74=\"\n$t.Cell(13, 4).Range.Text = \"74\u00f78=\"\n$t.Cell(13, 5).Range.Text = \"56\u00f76=\"\n\n$t.Cell(17, 1).Range.Text = \"18\u00f73=\"\n$t.Cell(17, 2).Range.Text = \"86\u00f79=\"\n$t.Cell(17, 3).Range.Text = \"99\u00f78=\"\n$t.Cell(17, 4).Range.Text = \"58\u00f79=\"\n$t.Cell(17, 5).Range.Text = \"92\u00f79=\"\n"}
